$wb = $excel.ActiveWorkbook

# --- system_permissions: insert a new permission row (SystemPermission / copy) ---
$ws1 = $wb.Worksheets.Item("system_permissions")
[void]$ws1.Activate()

# Insert a new row above the existing row 83 ("SystemPermission" / "share"),
# shifting it (and everything below) down by one.
[void]$ws1.Rows.Item(83).Insert()

$ws1.Range("A83").Value = "SystemPermission"
$ws1.Range("B83").Value = "copy"
$ws1.Range("C83").Value = "System Management"

# Match the C column formatting used by every other data row (style carried
# from the inserted row is already correct, but make sure explicitly).
[void]$ws1.Range("C82").Copy()
[void]$ws1.Range("C83").PasteSpecial(-4122) # xlPasteFormats

# --- user_assignments: it is no longer the active/selected tab ---
$ws4 = $wb.Worksheets.Item("user_assignments")
[void]$ws4.Range("A8").Select()

# system_permissions becomes (and remains) the active sheet/tab, with the
# selection landing on the row that used to hold "share" (now pushed to B88).
[void]$ws1.Activate()
[void]$ws1.Range("B88").Select()
